# Refactor synthetic array: insert a new "statut_name" column after
# statut_label (B) and before NCTId (old C, now D), shifting the existing
# C:L columns to D:M. Populate the new column with a human readable label
# derived from the existing statut_label (B) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C — shifts old C:L (NCTId..intervention_type) to D:M
$ws.Columns("C:C").Insert()

# New header
$ws.Range("C1").Value = "statut_name"

# Map statut_label -> statut_name for each data row (2-20)
$labels = @{
    "rouge"  = "résultat et / ou publication posté"
    "noir"   = "pas de résultat ni de publication"
    "vert"   = "résultat et / ou publication posté dans les 12 mois"
    "orange" = "résultat et / ou publication posté dans les 36 mois"
}

for ($r = 2; $r -le 20; $r++) {
    $label = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 3).Value = $labels[[string]$label]
}
